$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") for rows 2 through 484 all currently hold the
# serial date value 45203 (2023-10-04). Update them all to 45205
# (2023-10-06), matching the author's commit.
$ws.Range("C2:C484").Value = 45205
